$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.031.98'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '1.821.25'
$ws.Range('E3').Value = '  +3.24%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.89%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.04'
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4316'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3705'
$ws.Range('E8').Value = '  +2.70%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07289'
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').Value = '2.142.69'
$ws.Range('E10').Value = '  +21.84%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8720'
$ws.Range('E11').Value = '  +4.66%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.35'
$ws.Range('E12').Value = '  +5.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.425'
$ws.Range('E13').Value = '  +3.73%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.654'
$ws.Range('E14').Value = '  +3.83%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06995'
$ws.Range('E15').Value = '  +3.14%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '81.23'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.017'
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008943'
$ws.Range('E18').Value = '  +3.62%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.009'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.35'
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('D21').Value = '27.096.72'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.225'
$ws.Range('E22').Value = '  +4.53%  '
$ws.Range('D23').Value = '2.417.68'
$ws.Range('E23').Value = '  +21.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.04'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.92'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.45'
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.250'
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.928'
$ws.Range('E29').Value = '  +14.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.07'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08985'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.177'
$ws.Range('E32').Value = '  +6.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7496'
$ws.Range('E33').Value = '  +3.55%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.443'
$ws.Range('E34').Value = '  +2.95%  '
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.010'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  +5.55%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05251'
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01932'
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5148'
$ws.Range('E40').Value = '  +5.08%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1657'
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.748'
$ws.Range('E42').Value = '  +9.79%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.518'
$ws.Range('E43').Value = '  +5.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.360'
$ws.Range('E44').Value = '  +4.32%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '107.47'
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.42'
$ws.Range('E46').Value = '  +3.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.010'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4603'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.659'
$ws.Range('E49').Value = '  +5.58%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06325'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.824'
$ws.Range('E51').Value = '  +6.31%  '
